$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The steel description text had "RME/" removed from the 4th bullet
# (15% S/LFM+CDN/RME/H:1 -> 15% S/LFM+CDN/H:1). Rewrite the full
# multi-line cell value with the corrected text.
$ws.Range("B2").Value = "5% CR/LFM+CDN/H:2`n35% CR+PC/LFM+CDN/H:1`n25% S+SL/LFM+CDN/H:1`n15% S/LFM+CDN/H:1`n20% W/LWAL+CDN/H:1"

# Turn on wrap text for the (now shorter) multi-line description so it
# still displays on multiple lines within the cell.
$ws.Range("B2").WrapText = $true

# Grow row 2 so the wrapped text is fully visible.
$ws.Rows.Item(2).RowHeight = 256

# Restore the on-screen selection over the description block.
$ws.Range("B2:B9").Select()
